# Update row 6 of column A: "kink (length)" -> "length kindedd"
# (the old "kink (length)" label is dropped; "length kindedd" replaces it,
#  "kink (percentage)" in column C is untouched in meaning)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "length kindedd"

# Move the active selection, matching the author's last cursor position
$ws.Range("E13").Select()
